$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Inscritos 502 -> 503
$ws.Range("E10").Value = 503

# Row 23: Inscritos 196 -> 197
$ws.Range("E23").Value = 197

# Row 25: Inscritos 255 -> 256, Pagos 124 -> 125, Inscricoes homologadas 124 -> 125
$ws.Range("E25").Value = 256
$ws.Range("F25").Value = 125
$ws.Range("H25").Value = 125

# Row 30: Inscritos 195 -> 196, Pagos 117 -> 118, Inscricoes homologadas 117 -> 118
$ws.Range("E30").Value = 196
$ws.Range("F30").Value = 118
$ws.Range("H30").Value = 118

# Row 35: Pagos 89 -> 90, Inscricoes homologadas 89 -> 90
$ws.Range("F35").Value = 90
$ws.Range("H35").Value = 90

# Row 45: Inscritos 136 -> 137, Pagos 68 -> 69, Inscricoes homologadas 68 -> 69
$ws.Range("E45").Value = 137
$ws.Range("F45").Value = 69
$ws.Range("H45").Value = 69

# Row 47: Inscritos 435 -> 436, Pagos 215 -> 216, Inscricoes homologadas 215 -> 216
$ws.Range("E47").Value = 436
$ws.Range("F47").Value = 216
$ws.Range("H47").Value = 216

# Row 48: Inscritos 197 -> 198
$ws.Range("E48").Value = 198

# Row 51: Inscritos 224 -> 225
$ws.Range("E51").Value = 225
